$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Move Robot32 to location (2, 9) and remove the toolkit."

$ws.Range("A3").Value = "Move Robot6 to location (5, 4) and remove the liquid spill."

$ws.Range("A4").Value = "Move Robot29 to location (5, 12) and remove the large debris."
$ws.Range("B4").Value = $false

$ws.Range("A5").Value = "Move Robot48 to location (6, 6) and remove the dust."
$ws.Range("B5").Value = $true

$ws.Range("A6").Value = "Move Robot41 to location (1, 8) and remove the grass."

$ws.Range("A7").Value = "Move Robot10 to location (9, 5) and remove the small debris."

$ws.Range("A8").Value = "Move Robot13 to location (10, 10) and remove the vehicle."

$ws.Range("A9").Value = "Move Robot23 to location (8, 2) and remove the construction materials."

$ws.Range("A10").Value = "Move Robot24 to location (11, 8) and remove the tree branches."

$ws.Range("A11").Value = "Move Robot15 to location (3, 8) and remove the screws."
